# Auto-generated edit script: updates LevePriceNQ/HQ-derived profit columns (H:N)
# across all eight job sheets, per the scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7307.9165
$ws.Range("I2").Value = 5962.5
$ws.Range("J2").Value = 9998.75
$ws.Range("K2").Value = 5962.5
$ws.Range("L2").Value = 9998.75
$ws.Range("M2").Value = -5849.5
$ws.Range("N2").Value = -10224.75
$ws.Range("H4").Value = 78240.234
$ws.Range("I4").Value = 112713.664
$ws.Range("J4").Value = 675
$ws.Range("K4").Value = 112713.664
$ws.Range("L4").Value = 675
$ws.Range("M4").Value = -112599.664
$ws.Range("N4").Value = -903
$ws.Range("H13").Value = 6925
$ws.Range("J13").Value = 6925
$ws.Range("L13").Value = 6925
$ws.Range("N13").Value = -7263
$ws.Range("H31").Value = 1399
$ws.Range("I31").Value = 1399
$ws.Range("K31").Value = 4197
$ws.Range("M31").Value = -3967
$ws.Range("H33").Value = 8208722.5
$ws.Range("I33").Value = 10259353
$ws.Range("J33").Value = 6200.25
$ws.Range("K33").Value = 10259353
$ws.Range("L33").Value = 6200.25
$ws.Range("M33").Value = -10259124
$ws.Range("N33").Value = -6658.25
$ws.Range("H39").Value = 1631.5294
$ws.Range("I39").Value = 1137.5555
$ws.Range("K39").Value = 3412.6665
$ws.Range("M39").Value = -3116.6665
$ws.Range("H42").Value = 921.2857
$ws.Range("I42").Value = 197.33333
$ws.Range("K42").Value = 591.99999
$ws.Range("M42").Value = -361.99999
$ws.Range("H51").Value = 50103628
$ws.Range("I51").Value = 205199.6
$ws.Range("J51").Value = 100002056
$ws.Range("K51").Value = 205199.6
$ws.Range("L51").Value = 100002056
$ws.Range("M51").Value = -204715.6
$ws.Range("N51").Value = -100003024
$ws.Range("H53").Value = 55555904
$ws.Range("I53").Value = 158.25
$ws.Range("J53").Value = 100000504
$ws.Range("K53").Value = 158.25
$ws.Range("L53").Value = 100000504
$ws.Range("M53").Value = 478.75
$ws.Range("N53").Value = -100001778
$ws.Range("H70").Value = 1521.8667
$ws.Range("J70").Value = 1666.5
$ws.Range("L70").Value = 4999.5
$ws.Range("N70").Value = -5539.5
$ws.Range("H73").Value = 1521.8667
$ws.Range("J73").Value = 1666.5
$ws.Range("L73").Value = 4999.5
$ws.Range("N73").Value = -6871.5
$ws.Range("H74").Value = 3619.9167
$ws.Range("I74").Value = 3619.9167
$ws.Range("K74").Value = 3619.9167
$ws.Range("M74").Value = -2683.9167
$ws.Range("H77").Value = 3619.9167
$ws.Range("I77").Value = 3619.9167
$ws.Range("K77").Value = 18099.5835
$ws.Range("M77").Value = -13419.5835
$ws.Range("H80").Value = 41667344
$ws.Range("J80").Value = 13158759
$ws.Range("L80").Value = 39476277
$ws.Range("N80").Value = -39478273
$ws.Range("H83").Value = 41667344
$ws.Range("J83").Value = 13158759
$ws.Range("L83").Value = 118428831
$ws.Range("N83").Value = -118438815
$ws.Range("H100").Value = 2856.5715
$ws.Range("I100").Value = 2699.8
$ws.Range("J100").Value = 3248.5
$ws.Range("K100").Value = 2699.8
$ws.Range("L100").Value = 3248.5
$ws.Range("M100").Value = -2158.8
$ws.Range("N100").Value = -4330.5
$ws.Range("H103").Value = 407.41666
$ws.Range("J103").Value = 559.6
$ws.Range("L103").Value = 1678.8
$ws.Range("N103").Value = -2850.8
$ws.Range("H113").Value = 3382.6
$ws.Range("I113").Value = 3229
$ws.Range("K113").Value = 3229
$ws.Range("M113").Value = 25
$ws.Range("H126").Value = 134500
$ws.Range("J126").Value = 134500
$ws.Range("L126").Value = 134500
$ws.Range("N126").Value = -144380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 15858.75
$ws.Range("I28").Value = 15858.75
$ws.Range("K28").Value = 15858.75
$ws.Range("M28").Value = -15666.75
$ws.Range("H35").Value = 5765.857
$ws.Range("I35").Value = 5745.5
$ws.Range("J35").Value = 5888
$ws.Range("K35").Value = 5745.5
$ws.Range("L35").Value = 5888
$ws.Range("M35").Value = -5339.5
$ws.Range("N35").Value = -6700
$ws.Range("H45").Value = 36397.83
$ws.Range("I45").Value = 41805.4
$ws.Range("K45").Value = 41805.4
$ws.Range("M45").Value = -41428.4
$ws.Range("H99").Value = 15858.75
$ws.Range("I99").Value = 15858.75
$ws.Range("K99").Value = 15858.75
$ws.Range("M99").Value = -12863.75
$ws.Range("H102").Value = 1529.1111
$ws.Range("I102").Value = 1522.7059
$ws.Range("J102").Value = 1638
$ws.Range("K102").Value = 1522.7059
$ws.Range("L102").Value = 1638
$ws.Range("M102").Value = 99.29410000000007
$ws.Range("N102").Value = -4882
$ws.Range("H118").Value = 69999
$ws.Range("J118").Value = 69999
$ws.Range("L118").Value = 69999
$ws.Range("N118").Value = -73313

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 2038.6
$ws.Range("I12").Value = 2038.6
$ws.Range("K12").Value = 2038.6
$ws.Range("M12").Value = -1870.6
$ws.Range("H86").Value = 2106.6875
$ws.Range("I86").Value = 1201.9166
$ws.Range("J86").Value = 4821
$ws.Range("K86").Value = 1201.9166
$ws.Range("L86").Value = 4821
$ws.Range("M86").Value = -78.91660000000002
$ws.Range("N86").Value = -7067
$ws.Range("H89").Value = 2106.6875
$ws.Range("I89").Value = 1201.9166
$ws.Range("J89").Value = 4821
$ws.Range("K89").Value = 6009.583000000001
$ws.Range("L89").Value = 24105
$ws.Range("M89").Value = -393.5830000000005
$ws.Range("N89").Value = -35337
$ws.Range("H94").Value = 1971.5428
$ws.Range("I94").Value = 1470.5
$ws.Range("J94").Value = 3419
$ws.Range("K94").Value = 1470.5
$ws.Range("L94").Value = 3419
$ws.Range("M94").Value = -1019.5
$ws.Range("N94").Value = -4321
$ws.Range("H99").Value = 10409.786
$ws.Range("I99").Value = 14804.223
$ws.Range("K99").Value = 14804.223
$ws.Range("M99").Value = -13306.223
$ws.Range("H105").Value = 9993.929
$ws.Range("I105").Value = 13752.25
$ws.Range("J105").Value = 4982.8335
$ws.Range("K105").Value = 13752.25
$ws.Range("L105").Value = 4982.8335
$ws.Range("M105").Value = -12005.25
$ws.Range("N105").Value = -8476.833500000001
$ws.Range("H107").Value = 8220.8125
$ws.Range("I107").Value = 10179.292
$ws.Range("K107").Value = 10179.292
$ws.Range("M107").Value = -8259.291999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 127.416664
$ws.Range("I7").Value = 136.72728
$ws.Range("K7").Value = 136.72728
$ws.Range("M7").Value = -23.72728000000001
$ws.Range("H62").Value = 9560.9
$ws.Range("I62").Value = 9567.333000000001
$ws.Range("J62").Value = 9551.25
$ws.Range("K62").Value = 9567.333000000001
$ws.Range("L62").Value = 9551.25
$ws.Range("M62").Value = -8943.333000000001
$ws.Range("N62").Value = -10799.25
$ws.Range("H65").Value = 9560.9
$ws.Range("I65").Value = 9567.333000000001
$ws.Range("J65").Value = 9551.25
$ws.Range("K65").Value = 47836.665
$ws.Range("L65").Value = 47756.25
$ws.Range("M65").Value = -44716.665
$ws.Range("N65").Value = -53996.25
$ws.Range("H86").Value = 10758
$ws.Range("I86").Value = 4733.1113
$ws.Range("K86").Value = 4733.1113
$ws.Range("M86").Value = -3610.1113
$ws.Range("H89").Value = 10758
$ws.Range("I89").Value = 4733.1113
$ws.Range("K89").Value = 23665.5565
$ws.Range("M89").Value = -18049.5565
$ws.Range("H106").Value = 51887
$ws.Range("J106").Value = 51887
$ws.Range("L106").Value = 51887
$ws.Range("N106").Value = -54411
$ws.Range("H137").Value = 88779.39999999999
$ws.Range("J137").Value = 88779.39999999999
$ws.Range("L137").Value = 88779.39999999999
$ws.Range("N137").Value = -98979.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H47").Value = 2366.6667
$ws.Range("J47").Value = 2366.6667
$ws.Range("L47").Value = 7100.000100000001
$ws.Range("N47").Value = -7962.000100000001
$ws.Range("H117").Value = 66676640
$ws.Range("I117").Value = 7000
$ws.Range("J117").Value = 83344050
$ws.Range("K117").Value = 21000
$ws.Range("L117").Value = 250032150
$ws.Range("M117").Value = -17558
$ws.Range("N117").Value = -250039034
$ws.Range("H119").Value = 111127000
$ws.Range("I119").Value = 166678000
$ws.Range("K119").Value = 500034000
$ws.Range("M119").Value = -500029162
$ws.Range("H120").Value = 12853
$ws.Range("I120").Value = 7029.5
$ws.Range("K120").Value = 21088.5
$ws.Range("M120").Value = -16250.5
$ws.Range("H131").Value = 50123.906
$ws.Range("I131").Value = 1869.4166
$ws.Range("K131").Value = 5608.2498
$ws.Range("M131").Value = -568.2497999999996
$ws.Range("H138").Value = 4532
$ws.Range("I138").Value = 4532
$ws.Range("K138").Value = 13596
$ws.Range("M138").Value = -8456
$ws.Range("H140").Value = 2655.7058
$ws.Range("I140").Value = 2126.7693
$ws.Range("K140").Value = 6380.3079
$ws.Range("M140").Value = -1200.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8992.5
$ws.Range("J10").Value = 8992.5
$ws.Range("L10").Value = 8992.5
$ws.Range("N10").Value = -9330.5
$ws.Range("H22").Value = 4990
$ws.Range("J22").Value = 4990
$ws.Range("L22").Value = 4990
$ws.Range("N22").Value = -6048
$ws.Range("H36").Value = 4344.143
$ws.Range("J36").Value = 4741.8
$ws.Range("L36").Value = 4741.8
$ws.Range("N36").Value = -5711.8
$ws.Range("H70").Value = 9633.333000000001
$ws.Range("I70").Value = 9950
$ws.Range("K70").Value = 9950
$ws.Range("M70").Value = -9680
$ws.Range("H73").Value = 9633.333000000001
$ws.Range("I73").Value = 9950
$ws.Range("K73").Value = 9950
$ws.Range("M73").Value = -9014
$ws.Range("H80").Value = 10152711
$ws.Range("J80").Value = 38597680
$ws.Range("L80").Value = 38597680
$ws.Range("N80").Value = -38599676
$ws.Range("H83").Value = 10152711
$ws.Range("J83").Value = 38597680
$ws.Range("L83").Value = 192988400
$ws.Range("N83").Value = -192998384
$ws.Range("H97").Value = 557.74286
$ws.Range("I97").Value = 572.25
$ws.Range("J97").Value = 526.0909
$ws.Range("K97").Value = 572.25
$ws.Range("L97").Value = 526.0909
$ws.Range("M97").Value = -76.25
$ws.Range("N97").Value = -1518.0909
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884
$ws.Range("H134").Value = 65000.43
$ws.Range("J134").Value = 65000.43
$ws.Range("L134").Value = 195001.29
$ws.Range("N134").Value = -200071.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 742.4474
$ws.Range("I16").Value = 618.3929000000001
$ws.Range("J16").Value = 1089.8
$ws.Range("K16").Value = 618.3929000000001
$ws.Range("L16").Value = 1089.8
$ws.Range("M16").Value = -448.3929000000001
$ws.Range("N16").Value = -1429.8
$ws.Range("H46").Value = 16239.7
$ws.Range("I46").Value = 29249.25
$ws.Range("J46").Value = 7566.6665
$ws.Range("K46").Value = 29249.25
$ws.Range("L46").Value = 7566.6665
$ws.Range("M46").Value = -29061.25
$ws.Range("N46").Value = -7942.6665
$ws.Range("H61").Value = 2933.3076
$ws.Range("I61").Value = 2933.3076
$ws.Range("K61").Value = 2933.3076
$ws.Range("M61").Value = -2731.3076
$ws.Range("H76").Value = 9984
$ws.Range("J76").Value = 9984
$ws.Range("L76").Value = 9984
$ws.Range("N76").Value = -10660
$ws.Range("H79").Value = 9984
$ws.Range("J79").Value = 9984
$ws.Range("L79").Value = 9984
$ws.Range("N79").Value = -12324
$ws.Range("H101").Value = 22757.143
$ws.Range("J101").Value = 22757.143
$ws.Range("L101").Value = 22757.143
$ws.Range("N101").Value = -29247.143
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 26632.666
$ws.Range("J103").Value = 26632.666
$ws.Range("L103").Value = 26632.666
$ws.Range("N103").Value = -28976.666
$ws.Range("H113").Value = 2933.3076
$ws.Range("I113").Value = 2933.3076
$ws.Range("K113").Value = 2933.3076
$ws.Range("M113").Value = -763.3076000000001
$ws.Range("H121").Value = 41000
$ws.Range("J121").Value = 41000
$ws.Range("L121").Value = 41000
$ws.Range("N121").Value = -44494
$ws.Range("H136").Value = 2144.0688
$ws.Range("I136").Value = 2024.9412
$ws.Range("J136").Value = 2312.8333
$ws.Range("K136").Value = 6074.8236
$ws.Range("L136").Value = 6938.499899999999
$ws.Range("M136").Value = -3524.8236
$ws.Range("N136").Value = -12038.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 342333
$ws.Range("J3").Value = 512499.5
$ws.Range("L3").Value = 512499.5
$ws.Range("N3").Value = -512727.5
$ws.Range("H4").Value = 10694.571
$ws.Range("J4").Value = 6215
$ws.Range("L4").Value = 6215
$ws.Range("N4").Value = -6441
$ws.Range("H14").Value = 16582.928
$ws.Range("J14").Value = 18488.908
$ws.Range("L14").Value = 18488.908
$ws.Range("N14").Value = -18824.908
$ws.Range("H17").Value = 4180.6
$ws.Range("I17").Value = 4180.6
$ws.Range("K17").Value = 4180.6
$ws.Range("M17").Value = -4008.6
$ws.Range("H19").Value = 2500
$ws.Range("I19").Value = 2500
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -2326
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 11489.75
$ws.Range("J23").Value = 11489.75
$ws.Range("L23").Value = 11489.75
$ws.Range("N23").Value = -11947.75
$ws.Range("H100").Value = 1664.625
$ws.Range("I100").Value = 1616.7142
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 3233.4284
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2692.4284
$ws.Range("N100").Value = -5082
$ws.Range("H109").Value = 41326.668
$ws.Range("J109").Value = 40995
$ws.Range("L109").Value = 40995
$ws.Range("N109").Value = -43769
$ws.Range("H124").Value = 84000
$ws.Range("J124").Value = 84000
$ws.Range("L124").Value = 84000
$ws.Range("N124").Value = -93820

